$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (target) table data for rows 2-18 (columns A=Name, B=Position, C=Team)
$names = @(
    "Donovan Mitchell",
    "Malik Beasley",
    "Tim Hardaway Jr.",
    "Alperen Sengün",
    "Domantas Sabonis",
    "Santi Aldama",
    "Kristaps Porzingis",
    "Michael Porter Jr.",
    "Chris Boucher",
    "Andrew Wiggins",
    "Kelly Oubre Jr.",
    "Dyson Daniels",
    "Josh Hart",
    "Victor Wembanyama",
    "Donte DiVincenzo",
    "Cam Thomas",
    "De'Andre Hunter"
)

$positions = @(
    "PG,SG",
    "SG,SF",
    "SG,SF",
    "C",
    "C",
    "PF,C",
    "PF,C",
    "SF,PF",
    "PF,C",
    "SF,PF",
    "SG,SF",
    "PG,SG,SF",
    "SG,SF,PF",
    "C",
    "PG,SG,SF",
    "SG,SF",
    "SF,PF"
)

$teams = @(
    "Cleveland Cavaliers",
    "Detroit Pistons",
    "Detroit Pistons",
    "Houston Rockets",
    "Sacramento Kings",
    "Memphis Grizzlies",
    "Boston Celtics",
    "Denver Nuggets",
    "Toronto Raptors",
    "Golden State Warriors",
    "Philadelphia 76ers",
    "Atlanta Hawks",
    "New York Knicks",
    "San Antonio Spurs",
    "Minnesota Timberwolves",
    "Brooklyn Nets",
    "Atlanta Hawks"
)

# Clear out the previous data rows before rewriting.
$ws.Range("A2:C18").Clear()

$row = 2
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $row++
}

$row = 2
foreach ($pos in $positions) {
    $ws.Cells.Item($row, 2).Value = $pos
    $row++
}

$row = 2
foreach ($team in $teams) {
    $ws.Cells.Item($row, 3).Value = $team
    $row++
}
